# Remove the trailing "Ver no Jupiter..." / copyright footer block, plus
# the blank paragraph that precedes it, leaving the final blank paragraph
# and the page-break paragraph untouched.

$d = $word.ActiveDocument

$startMarker = "Ver no Jupiter"
$endMarker = "Powered by Jekyll"

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -like "*$startMarker*") {
        $startIdx = $i
    }
    if ($t -like "*$endMarker*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1) {
    # Also drop the blank paragraph immediately before the "Ver no
    # Jupiter..." paragraph, matching the diff.
    $delStartIdx = $startIdx
    if ($startIdx -gt 1 -and $d.Paragraphs.Item($startIdx - 1).Range.Text -eq "`r") {
        $delStartIdx = $startIdx - 1
    }

    $delStart = $d.Paragraphs.Item($delStartIdx).Range.Start
    $delEnd = $d.Paragraphs.Item($endIdx).Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
